$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1039.5416
$ws.Range("I132").Value = 717.119
$ws.Range("K132").Value = 2151.357
$ws.Range("M132").Value = 378.643
$ws.Range("H135").Value = 1250897.5
$ws.Range("I135").Value = 2222706.5
$ws.Range("J135").Value = 1429
$ws.Range("K135").Value = 20004358.5
$ws.Range("L135").Value = 12861
$ws.Range("M135").Value = -20001823.5
$ws.Range("N135").Value = -17931
$ws.Range("H137").Value = 4114.826
$ws.Range("J137").Value = 7474.4707
$ws.Range("L137").Value = 22423.4121
$ws.Range("N137").Value = -27523.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3402.7693
$ws.Range("I2").Value = 3729.7144
$ws.Range("J2").Value = 3021.3333
$ws.Range("K2").Value = 3729.7144
$ws.Range("L2").Value = 3021.3333
$ws.Range("M2").Value = -3616.7144
$ws.Range("N2").Value = -3247.3333
$ws.Range("H32").Value = 3181355.8
$ws.Range("I32").Value = 3453034.8
$ws.Range("K32").Value = 3453034.8
$ws.Range("M32").Value = -3452747.8
$ws.Range("H45").Value = 2945.9565
$ws.Range("I45").Value = 2131.8
$ws.Range("J45").Value = 3572.2307
$ws.Range("K45").Value = 2131.8
$ws.Range("L45").Value = 3572.2307
$ws.Range("M45").Value = -1754.8
$ws.Range("N45").Value = -4326.2307
$ws.Range("H61").Value = 22729592
$ws.Range("I61").Value = 1719.0286
$ws.Range("J61").Value = 111115770
$ws.Range("K61").Value = 1719.0286
$ws.Range("L61").Value = 111115770
$ws.Range("M61").Value = -1507.0286
$ws.Range("N61").Value = -111116194
$ws.Range("H110").Value = 14498275
$ws.Range("J110").Value = 66667252
$ws.Range("L110").Value = 66667252
$ws.Range("N110").Value = -66671342
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H116").Value = 3402.7693
$ws.Range("I116").Value = 3729.7144
$ws.Range("J116").Value = 3021.3333
$ws.Range("K116").Value = 3729.7144
$ws.Range("L116").Value = 3021.3333
$ws.Range("M116").Value = -1435.7144
$ws.Range("N116").Value = -7609.3333
$ws.Range("H122").Value = 5438.278
$ws.Range("I122").Value = 3293.8572
$ws.Range("K122").Value = 9881.571599999999
$ws.Range("M122").Value = -7431.571599999999
$ws.Range("H132").Value = 5610.1816
$ws.Range("I132").Value = 5307
$ws.Range("J132").Value = 6259.857
$ws.Range("K132").Value = 15921
$ws.Range("L132").Value = 18779.571
$ws.Range("M132").Value = -13391
$ws.Range("N132").Value = -23839.571
$ws.Range("H136").Value = 22729592
$ws.Range("I136").Value = 1719.0286
$ws.Range("J136").Value = 111115770
$ws.Range("K136").Value = 5157.085800000001
$ws.Range("L136").Value = 333347310
$ws.Range("M136").Value = -2607.085800000001
$ws.Range("N136").Value = -333352410

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3402.7693
$ws.Range("I3").Value = 3729.7144
$ws.Range("J3").Value = 3021.3333
$ws.Range("K3").Value = 3729.7144
$ws.Range("L3").Value = 3021.3333
$ws.Range("M3").Value = -3615.7144
$ws.Range("N3").Value = -3249.3333
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H86").Value = 45560.26
$ws.Range("I86").Value = 72442.86
$ws.Range("J86").Value = 3742.889
$ws.Range("K86").Value = 72442.86
$ws.Range("L86").Value = 3742.889
$ws.Range("M86").Value = -71319.86
$ws.Range("N86").Value = -5988.889
$ws.Range("H89").Value = 45560.26
$ws.Range("I89").Value = 72442.86
$ws.Range("J89").Value = 3742.889
$ws.Range("K89").Value = 362214.3
$ws.Range("L89").Value = 18714.445
$ws.Range("M89").Value = -356598.3
$ws.Range("N89").Value = -29946.445
$ws.Range("H105").Value = 2581.9363
$ws.Range("I105").Value = 1657.2059
$ws.Range("K105").Value = 1657.2059
$ws.Range("M105").Value = 89.79410000000007

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4315.477
$ws.Range("I31").Value = 1756.9131
$ws.Range("J31").Value = 7117.7144
$ws.Range("K31").Value = 1756.9131
$ws.Range("L31").Value = 7117.7144
$ws.Range("M31").Value = -1461.9131
$ws.Range("N31").Value = -7707.7144
$ws.Range("H34").Value = 4315.477
$ws.Range("I34").Value = 1756.9131
$ws.Range("J34").Value = 7117.7144
$ws.Range("K34").Value = 1756.9131
$ws.Range("L34").Value = 7117.7144
$ws.Range("M34").Value = -1554.9131
$ws.Range("N34").Value = -7521.7144
$ws.Range("H99").Value = 7426.625
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 7426.625
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H134").Value = 3768.282
$ws.Range("I134").Value = 2065.6316
$ws.Range("K134").Value = 6196.8948
$ws.Range("M134").Value = -3661.8948
$ws.Range("H138").Value = 44500
$ws.Range("J138").Value = 44500
$ws.Range("L138").Value = 44500
$ws.Range("N138").Value = -54780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5882.4375
$ws.Range("I132").Value = 2624.5
$ws.Range("J132").Value = 7837.2
$ws.Range("K132").Value = 23620.5
$ws.Range("L132").Value = 70534.8
$ws.Range("M132").Value = -21090.5
$ws.Range("N132").Value = -75594.8
$ws.Range("H133").Value = 13163.77
$ws.Range("I133").Value = 9018.429
$ws.Range("J133").Value = 18000
$ws.Range("K133").Value = 27055.287
$ws.Range("L133").Value = 54000
$ws.Range("M133").Value = -21995.287
$ws.Range("N133").Value = -64120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1054.2106
$ws.Range("I97").Value = 921.25
$ws.Range("J97").Value = 1763.3334
$ws.Range("K97").Value = 921.25
$ws.Range("L97").Value = 1763.3334
$ws.Range("M97").Value = -425.25
$ws.Range("N97").Value = -2755.3334
$ws.Range("H122").Value = 3634227.2
$ws.Range("I122").Value = 3825291.8
$ws.Range("K122").Value = 11475875.4
$ws.Range("M122").Value = -11473425.4
$ws.Range("H132").Value = 1964.0731
$ws.Range("I132").Value = 1306.875
$ws.Range("J132").Value = 4300.778
$ws.Range("K132").Value = 3920.625
$ws.Range("L132").Value = 12902.334
$ws.Range("M132").Value = -1390.625
$ws.Range("N132").Value = -17962.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4895.6
$ws.Range("I7").Value = 4128.75
$ws.Range("J7").Value = 5772
$ws.Range("K7").Value = 4128.75
$ws.Range("L7").Value = 5772
$ws.Range("M7").Value = -4016.75
$ws.Range("N7").Value = -5996
$ws.Range("H22").Value = 1233.7826
$ws.Range("I22").Value = 756.1579
$ws.Range("K22").Value = 756.1579
$ws.Range("M22").Value = -461.1579
$ws.Range("H27").Value = 1233.7826
$ws.Range("I27").Value = 756.1579
$ws.Range("K27").Value = 756.1579
$ws.Range("M27").Value = -649.1579
$ws.Range("H40").Value = 5249.55
$ws.Range("I40").Value = 2200.8
$ws.Range("K40").Value = 2200.8
$ws.Range("M40").Value = -2064.8
$ws.Range("H46").Value = 2147.6553
$ws.Range("I46").Value = 1513.6471
$ws.Range("J46").Value = 3045.8333
$ws.Range("K46").Value = 1513.6471
$ws.Range("L46").Value = 3045.8333
$ws.Range("M46").Value = -1325.6471
$ws.Range("N46").Value = -3421.8333
$ws.Range("H122").Value = 4331.391
$ws.Range("J122").Value = 6849.625
$ws.Range("L122").Value = 20548.875
$ws.Range("N122").Value = -25448.875
$ws.Range("H126").Value = 4895.6
$ws.Range("I126").Value = 4128.75
$ws.Range("J126").Value = 5772
$ws.Range("K126").Value = 12386.25
$ws.Range("L126").Value = 17316
$ws.Range("M126").Value = -9916.25
$ws.Range("N126").Value = -22256
$ws.Range("H132").Value = 4738.8535
$ws.Range("J132").Value = 6679.9
$ws.Range("L132").Value = 20039.7
$ws.Range("N132").Value = -25099.7
$ws.Range("H136").Value = 11162.639
$ws.Range("I136").Value = 3053.4119
$ws.Range("J136").Value = 18418.264
$ws.Range("K136").Value = 9160.235700000001
$ws.Range("L136").Value = 55254.792
$ws.Range("M136").Value = -6610.235700000001
$ws.Range("N136").Value = -60354.792

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 824.4039
$ws.Range("I113").Value = 733.8108
$ws.Range("J113").Value = 1047.8667
$ws.Range("K113").Value = 2201.4324
$ws.Range("L113").Value = 3143.6001
$ws.Range("M113").Value = -31.43239999999969
$ws.Range("N113").Value = -7483.6001
$ws.Range("H122").Value = 406300
$ws.Range("I122").Value = 1004200.75
$ws.Range("K122").Value = 3012602.25
$ws.Range("M122").Value = -3010152.25
$ws.Range("H132").Value = 5766.4136
$ws.Range("I132").Value = 6756
$ws.Range("K132").Value = 20268
$ws.Range("M132").Value = -17738
$ws.Range("H136").Value = 482683.9
$ws.Range("I136").Value = 2251
$ws.Range("J136").Value = 595726.9399999999
$ws.Range("K136").Value = 6753
$ws.Range("L136").Value = 1787180.82
$ws.Range("M136").Value = -4203
$ws.Range("N136").Value = -1792280.82
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
